# Updates crypto price/volume table cells to the latest scraped values.
# Numeric-looking price text (e.g. "602.43") is prefixed with a leading
# apostrophe so Excel stores it as text, matching the original inlineStr
# cells instead of converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.868.97"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "3.742.29"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'602.43"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "'168.67"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("D7").Value = "3.739.70"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("D11").Value = "'6.34"
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "'38.25"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D15").Value = "4.369.44"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "3.744.76"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "68.885.24"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "'10.84"
$ws.Range("E21").Value = "  +19.05%  "
$ws.Range("D22").Value = "'493.33"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "'0.726"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "'0.0000155"
$ws.Range("E24").Value = "  +12.88%  "
$ws.Range("D25").Value = "'85.23"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'2.32"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("D28").Value = "'10.39"
$ws.Range("E28").Value = "  +4.52%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  +7.01%  "
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("D32").Value = "'7.97"
$ws.Range("E32").Value = "  +3.96%  "
$ws.Range("D33").Value = "'31.74"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "3.889.31"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("D35").Value = "'0.108"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").Value = "3.678.12"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("D39").Value = "'5.85"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'3.00"
$ws.Range("E42").Value = "  +8.45%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'436.28"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").Value = "'48.87"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D48").Value = "'40.45"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "'141.48"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").Value = "'0.0355"
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("D51").Value = "2.771.90"
$ws.Range("E51").Value = "  +0.91%  "
